$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New forecast row for the latest date_of_forecast (2025-11-25), carrying
# over the same date-cell formatting (s="2") used by the rows above it.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -0.7200474048664085
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = -2.181280391105744
